$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 9
$ws.Range("H9").Value = 116.666664
$ws.Range("I9").Value = 183.33333
$ws.Range("J9").Value = 50
$ws.Range("K9").Value = 183.33333
$ws.Range("L9").Value = 50
$ws.Range("M9").Value = -14.33332999999999
$ws.Range("N9").Value = -388

# Row 51
$ws.Range("H51").Value = 3594.4443
$ws.Range("J51").Value = 3807.1428
$ws.Range("L51").Value = 3807.1428
$ws.Range("N51").Value = -4775.1428

# Row 58
$ws.Range("H58").Value = 953.9
$ws.Range("J58").Value = 3000
$ws.Range("L58").Value = 9000
$ws.Range("N58").Value = -9300

# Row 64
$ws.Range("H64").Value = 8295.200000000001
$ws.Range("I64").Value = 7326
$ws.Range("J64").Value = 9749
$ws.Range("K64").Value = 7326
$ws.Range("L64").Value = 9749
$ws.Range("M64").Value = -7078
$ws.Range("N64").Value = -10245

# Row 67
$ws.Range("H67").Value = 8295.200000000001
$ws.Range("I67").Value = 7326
$ws.Range("J67").Value = 9749
$ws.Range("K67").Value = 7326
$ws.Range("L67").Value = 9749
$ws.Range("M67").Value = -6468
$ws.Range("N67").Value = -11465

# Row 76
$ws.Range("H76").Value = 5000
$ws.Range("J76").Value = 2000
$ws.Range("L76").Value = 2000
$ws.Range("N76").Value = -2630

# Row 79
$ws.Range("H79").Value = 5000
$ws.Range("J79").Value = 2000
$ws.Range("L79").Value = 2000
$ws.Range("N79").Value = -4184

# Row 105
$ws.Range("H105").Value = 34835.5
$ws.Range("J105").Value = 34835.5
$ws.Range("L105").Value = 34835.5
$ws.Range("N105").Value = -41823.5

# Row 107
$ws.Range("H107").Value = 123.5
$ws.Range("I107").Value = 123.5
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 123.5
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1796.5
$ws.Range("N107").ClearContents()

# Row 132
$ws.Range("H132").Value = 13589.869
$ws.Range("I132").Value = 12753.777
$ws.Range("J132").Value = 16599.8
$ws.Range("K132").Value = 38261.331
$ws.Range("L132").Value = 49799.39999999999
$ws.Range("M132").Value = -35731.331
$ws.Range("N132").Value = -54859.39999999999

# Row 137
$ws.Range("H137").Value = 3083.25
$ws.Range("I137").Value = 2025
$ws.Range("J137").Value = 4000.4
$ws.Range("K137").Value = 6075
$ws.Range("L137").Value = 12001.2
$ws.Range("M137").Value = -3525
$ws.Range("N137").Value = -17101.2

# Row 141
$ws.Range("H141").Value = 2075.1667
$ws.Range("I141").Value = 1145
$ws.Range("J141").Value = 4865.6665
$ws.Range("K141").Value = 3435
$ws.Range("L141").Value = 14596.9995
$ws.Range("M141").Value = 1745
$ws.Range("N141").Value = -24956.9995

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 2727.08
$ws.Range("I45").Value = 2418.9
$ws.Range("J45").Value = 3959.8
$ws.Range("K45").Value = 2418.9
$ws.Range("L45").Value = 3959.8
$ws.Range("M45").Value = -2041.9
$ws.Range("N45").Value = -4713.8

# Row 61
$ws.Range("H61").Value = 5049.8335
$ws.Range("I61").Value = 3824.75
$ws.Range("K61").Value = 3824.75
$ws.Range("M61").Value = -3612.75

# Row 110
$ws.Range("H110").Value = 2771.625
$ws.Range("I110").Value = 834.8
$ws.Range("K110").Value = 834.8
$ws.Range("M110").Value = 1210.2

# Row 122
$ws.Range("H122").Value = 1426.0769
$ws.Range("I122").Value = 1426.0769
$ws.Range("K122").Value = 4278.2307
$ws.Range("M122").Value = -1828.2307

# Row 136
$ws.Range("H136").Value = 5049.8335
$ws.Range("I136").Value = 3824.75
$ws.Range("K136").Value = 11474.25
$ws.Range("M136").Value = -8924.25

$ws = $wb.Worksheets.Item("BSM")
# Row 132
$ws.Range("H132").Value = 131666.67
$ws.Range("J132").Value = 131666.67
$ws.Range("L132").Value = 131666.67
$ws.Range("N132").Value = -141786.67

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 3149.7144
$ws.Range("I22").Value = 524.5
$ws.Range("J22").Value = 4199.8
$ws.Range("K22").Value = 524.5
$ws.Range("L22").Value = 4199.8
$ws.Range("M22").Value = -174.5
$ws.Range("N22").Value = -4899.8

# Row 132
$ws.Range("H132").Value = 3408.2307
$ws.Range("I132").Value = 3048.111
$ws.Range("K132").Value = 9144.332999999999
$ws.Range("M132").Value = -6614.332999999999

$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 545957.9399999999
$ws.Range("I4").Value = 2600600.2
$ws.Range("J4").Value = 5262.6313
$ws.Range("K4").Value = 7801800.600000001
$ws.Range("L4").Value = 15787.8939
$ws.Range("M4").Value = -7801688.600000001
$ws.Range("N4").Value = -16011.8939

# Row 5
$ws.Range("H5").Value = 1585.25
$ws.Range("J5").Value = 2299
$ws.Range("L5").Value = 6897
$ws.Range("N5").Value = -7121

# Row 132
$ws.Range("H132").Value = 2180.8
$ws.Range("I132").Value = 2101
$ws.Range("K132").Value = 18909
$ws.Range("M132").Value = -16379

# Row 135
$ws.Range("H135").Value = 1585.25
$ws.Range("J135").Value = 2299
$ws.Range("L135").Value = 20691
$ws.Range("N135").Value = -25761

# Row 140
$ws.Range("H140").Value = 1888.3572
$ws.Range("J140").Value = 5449.5
$ws.Range("L140").Value = 16348.5
$ws.Range("N140").Value = -26708.5

$ws = $wb.Worksheets.Item("GSM")
# Row 61
$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()

# Row 80
$ws.Range("H80").Value = 3299.3333
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()

# Row 83
$ws.Range("H83").Value = 3299.3333
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

# Row 122
$ws.Range("H122").Value = 459354.38
$ws.Range("I122").Value = 560017.9
$ws.Range("K122").Value = 1680053.7
$ws.Range("M122").Value = -1677603.7

# Row 126
$ws.Range("H126").Value = 4025.4285
$ws.Range("I126").Value = 3779.6667
$ws.Range("K126").Value = 11339.0001
$ws.Range("M126").Value = -8869.000100000001

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 548.1667
$ws.Range("I16").Value = 557.8
$ws.Range("K16").Value = 557.8
$ws.Range("M16").Value = -387.8

# Row 46
$ws.Range("H46").Value = 8230.786
$ws.Range("I46").Value = 9308.25
$ws.Range("J46").Value = 7799.8
$ws.Range("K46").Value = 9308.25
$ws.Range("L46").Value = 7799.8
$ws.Range("M46").Value = -9120.25
$ws.Range("N46").Value = -8175.8

# Row 63
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()

# Row 66
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()

# Row 69
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()

# Row 72
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()

# Row 132
$ws.Range("H132").Value = 5853.857
$ws.Range("I132").Value = 5162
$ws.Range("K132").Value = 15486
$ws.Range("M132").Value = -12956

$ws = $wb.Worksheets.Item("WVR")
# Row 41
$ws.Range("H41").Value = 39938.5
$ws.Range("J41").Value = 39918
$ws.Range("L41").Value = 39918
$ws.Range("N41").Value = -40698

# Row 126
$ws.Range("H126").Value = 3953.3
$ws.Range("I126").Value = 1839.75
$ws.Range("J126").Value = 7123.625
$ws.Range("K126").Value = 5519.25
$ws.Range("L126").Value = 21370.875
$ws.Range("M126").Value = -3049.25
$ws.Range("N126").Value = -26310.875

# Row 132
$ws.Range("H132").Value = 4220
$ws.Range("I132").Value = 4016.0667
$ws.Range("J132").Value = 5749.5
$ws.Range("K132").Value = 12048.2001
$ws.Range("L132").Value = 17248.5
$ws.Range("M132").Value = -9518.2001
$ws.Range("N132").Value = -22308.5

# Row 136
$ws.Range("H136").Value = 4396.727
$ws.Range("I136").Value = 3546.75
$ws.Range("J136").Value = 6663.3335
$ws.Range("K136").Value = 10640.25
$ws.Range("L136").Value = 19990.0005
$ws.Range("M136").Value = -8090.25
$ws.Range("N136").Value = -25090.0005
